# Auto-generated: apply cell-value updates to Halicarnassus_Profits sheets
# (market data refresh — currentAveragePrice* / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 195.92857
$ws.Range("I39").Value = 137
$ws.Range("K39").Value = 411
$ws.Range("M39").Value = -115
$ws.Range("H88").Value = 2129.75
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 1839.6666
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 1839.6666
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -2651.6666
$ws.Range("H91").Value = 2129.75
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 1839.6666
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 1839.6666
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -4647.6666
$ws.Range("H111").Value = 500
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = $null
$ws.Range("H116").Value = 4998.75
$ws.Range("J116").Value = 5331.6665
$ws.Range("L116").Value = 5331.6665
$ws.Range("N116").Value = -12215.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1021.625
$ws.Range("I2").Value = 1009.73334
$ws.Range("K2").Value = 1009.73334
$ws.Range("M2").Value = -896.73334
$ws.Range("H6").Value = 10590938
$ws.Range("I6").Value = 8574711
$ws.Range("J6").Value = 20000000
$ws.Range("K6").Value = 8574711
$ws.Range("L6").Value = 20000000
$ws.Range("M6").Value = -8574538
$ws.Range("N6").Value = -20000346
$ws.Range("H45").Value = 2759.3333
$ws.Range("I45").Value = 1903.3636
$ws.Range("K45").Value = 1903.3636
$ws.Range("M45").Value = -1526.3636
$ws.Range("H116").Value = 1021.625
$ws.Range("I116").Value = 1009.73334
$ws.Range("K116").Value = 1009.73334
$ws.Range("M116").Value = 1284.26666
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1021.625
$ws.Range("I3").Value = 1009.73334
$ws.Range("K3").Value = 1009.73334
$ws.Range("M3").Value = -895.73334
$ws.Range("H86").Value = 2828.6191
$ws.Range("I86").Value = 1226.8
$ws.Range("K86").Value = 1226.8
$ws.Range("M86").Value = -103.8
$ws.Range("H89").Value = 2828.6191
$ws.Range("I89").Value = 1226.8
$ws.Range("K89").Value = 6134
$ws.Range("M89").Value = -518
$ws.Range("H96").Value = 18419.834
$ws.Range("I96").Value = 18419.834
$ws.Range("K96").Value = 18419.834
$ws.Range("M96").Value = -15673.834
$ws.Range("H112").Value = 600
$ws.Range("J112").Value = 600
$ws.Range("L112").Value = 600
$ws.Range("N112").Value = -3554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 43686.2
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 52107.75
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 52107.75
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -53579.75
$ws.Range("H60").Value = 57687
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 73582.664
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 73582.664
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -74604.664
$ws.Range("H61").Value = 43686.2
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 52107.75
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 52107.75
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -52803.75
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H99").Value = 3599.8
$ws.Range("J99").Value = 3499.5
$ws.Range("L99").Value = 3499.5
$ws.Range("N99").Value = -6495.5
$ws.Range("H126").Value = 3599.8
$ws.Range("J126").Value = 3499.5
$ws.Range("L126").Value = 10498.5
$ws.Range("N126").Value = -15438.5
$ws.Range("H132").Value = 1993.2222
$ws.Range("I132").Value = 2028.2941
$ws.Range("K132").Value = 6084.8823
$ws.Range("M132").Value = -3554.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1416.9375
$ws.Range("J5").Value = 1449.8334
$ws.Range("L5").Value = 4349.5002
$ws.Range("N5").Value = -4573.5002
$ws.Range("H128").Value = 555737
$ws.Range("I128").Value = 555737
$ws.Range("K128").Value = 1667211
$ws.Range("M128").Value = -1662231
$ws.Range("H135").Value = 1416.9375
$ws.Range("J135").Value = 1449.8334
$ws.Range("L135").Value = 13048.5006
$ws.Range("N135").Value = -18118.5006
$ws.Range("H138").Value = 6571.4614
$ws.Range("I138").Value = 1925.8
$ws.Range("K138").Value = 5777.4
$ws.Range("M138").Value = -637.3999999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 111666.336
$ws.Range("J133").Value = 111666.336
$ws.Range("L133").Value = 111666.336
$ws.Range("N133").Value = -121786.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1400
$ws.Range("I7").Value = 1437.5
$ws.Range("K7").Value = 1437.5
$ws.Range("M7").Value = -1325.5
$ws.Range("H22").Value = 1104.0625
$ws.Range("I22").Value = 694.6667
$ws.Range("J22").Value = 1349.7
$ws.Range("K22").Value = 694.6667
$ws.Range("L22").Value = 1349.7
$ws.Range("M22").Value = -399.6667
$ws.Range("N22").Value = -1939.7
$ws.Range("H27").Value = 1104.0625
$ws.Range("I27").Value = 694.6667
$ws.Range("J27").Value = 1349.7
$ws.Range("K27").Value = 694.6667
$ws.Range("L27").Value = 1349.7
$ws.Range("M27").Value = -587.6667
$ws.Range("N27").Value = -1563.7
$ws.Range("H40").Value = 2251.3333
$ws.Range("I40").Value = 2251.3333
$ws.Range("K40").Value = 2251.3333
$ws.Range("M40").Value = -2115.3333
$ws.Range("H42").Value = 20025
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = $null
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = $null
$ws.Range("H49").Value = 20025
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null
$ws.Range("H68").Value = 5930
$ws.Range("J68").Value = 6937.75
$ws.Range("L68").Value = 6937.75
$ws.Range("N68").Value = -8435.75
$ws.Range("H71").Value = 5930
$ws.Range("J71").Value = 6937.75
$ws.Range("L71").Value = 34688.75
$ws.Range("N71").Value = -42176.75
$ws.Range("H82").Value = 2611.389
$ws.Range("J82").Value = 4239.8
$ws.Range("L82").Value = 4239.8
$ws.Range("N82").Value = -4961.8
$ws.Range("H85").Value = 2611.389
$ws.Range("J85").Value = 4239.8
$ws.Range("L85").Value = 4239.8
$ws.Range("N85").Value = -6735.8
$ws.Range("H126").Value = 1400
$ws.Range("I126").Value = 1437.5
$ws.Range("K126").Value = 4312.5
$ws.Range("M126").Value = -1842.5
$ws.Range("H132").Value = 9499.5
$ws.Range("I132").Value = 9499.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 28498.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -25968.5
$ws.Range("N132").Value = $null
$ws.Range("H136").Value = 2562.25
$ws.Range("I136").Value = 2562.25
$ws.Range("K136").Value = 7686.75
$ws.Range("M136").Value = -5136.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = $null
$ws.Range("H32").Value = 4700.5
$ws.Range("J32").Value = 1400
$ws.Range("L32").Value = 1400
$ws.Range("N32").Value = -2034
$ws.Range("H100").Value = 795
$ws.Range("I100").Value = 795
$ws.Range("K100").Value = 1590
$ws.Range("M100").Value = -1049
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H136").Value = 3251.75
$ws.Range("I136").Value = 2015.875
$ws.Range("K136").Value = 6047.625
$ws.Range("M136").Value = -3497.625
